$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AVIV_CheckoutPage")

# D6 and D8 hold long numeric-looking identifiers that must stay stored as
# text (matching the original sharedString cell type). Force text entry by
# temporarily formatting as Text, then restore the default "Normal" style
# so the cell's style index is left untouched.
$d6 = $ws.Range("D6")
$d6.NumberFormat = "@"
$d6.Value = "3518152942"
$d6.Style = "Normal"

$d8 = $ws.Range("D8")
$d8.NumberFormat = "@"
$d8.Value = "5894143326"
$d8.Style = "Normal"

# D13 is plain text already, no special handling required.
$ws.Range("D13").Value = "ORDER NUMBER: 1034"
